# Add a new "Complaint" participant-assignment rule row (row 21) to Sheet1,
# mirroring the existing "Case File" rule row (row 20), and update the
# sheet's view selection accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy formatting (styles, number formats, etc.) from the source row's B:F
# cells onto the new row's B:F cells - mirrors the existing "Case File" rule
# row without touching column A (the diff doesn't emit an A21 cell).
$ws.Range("B20:F20").Copy($ws.Range("B21:F21"))
$ws.Rows.Item(21).RowHeight = $ws.Rows.Item(20).RowHeight

# Populate the new row's values.
$ws.Cells.Item(21, 2).Value2 = "Complaint - Check participants list for NoAccess & Owner"
$ws.Cells.Item(21, 3).Value2 = "COMPLAINT"
$ws.Cells.Item(21, 4).Value2 = $ws.Cells.Item(20, 4).Value2
$ws.Cells.Item(21, 5).Value2 = $ws.Cells.Item(20, 5).Value2
$ws.Cells.Item(21, 6).Value2 = $ws.Cells.Item(20, 6).Value2

# Update the visible view state to match the post-edit workbook.
$ws.Range("D21").Select()
